$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 62: Comments (column F) count 1 -> 3
$ws.Cells.Item(62, 6).Value = 3

# Append new row 63
$ws.Cells.Item(63, 1).Value = 44148
$ws.Cells.Item(63, 2).Value = 0.30902777777777779
$ws.Cells.Item(63, 3).Value = "Friends"
$ws.Cells.Item(63, 4).Value = "Dense Fog Friday ☁️"
$ws.Cells.Item(63, 5).Value = "10107821070853229"
$ws.Cells.Item(63, 6).Value = 3
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 8).Value = 0
$ws.Cells.Item(63, 9).Value = 2
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 11).Value = 3
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 13).Value = 0

# Match style/number-format of the date/time/text columns to the rest of the table
$ws.Cells.Item(63, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(63, 2).NumberFormat = "h:mm:ss;@"
$ws.Cells.Item(63, 5).NumberFormat = "@"

# Update the view state to match the saved workbook
$ws.Application.ActiveWindow.ScrollRow = 43
$ws.Range("E60").Select() | Out-Null
